$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 7259
$ws1.Range("F16").Value = 122
$ws1.Range("F18").Value = 11960
$ws1.Range("F21").Value = 2395
$ws1.Range("F23").Value = 3393
$ws1.Range("F26").Value = 2834
$ws1.Range("F33").Value = 2431
$ws1.Range("F35").Value = 1679
$ws1.Range("F38").Value = 5940
$ws1.Range("F40").Value = 22

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 255

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 294

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 294
$ws4.Range("F11").Value = 7259
$ws4.Range("F12").Value = 7259
$ws4.Range("F19").Value = 255
$ws4.Range("F23").Value = 11960
$ws4.Range("F26").Value = 2395
$ws4.Range("F27").Value = 2395
$ws4.Range("F28").Value = 3393
$ws4.Range("F29").Value = 2834
$ws4.Range("F36").Value = 2431
$ws4.Range("F38").Value = 1679
$ws4.Range("F40").Value = 5940
